$wb = $excel.ActiveWorkbook

# ----- Sheet1 -----
$ws1 = $wb.Worksheets.Item("Sheet1")

# Clear old data range (A2:B11) before writing the new, shorter table
$ws1.Range("A2:B11").ClearContents()

$sheet1Data = @(
    @(1, 6),
    @(2, 9),
    @(3, 8),
    @(141, 142),
    @(142, 143),
    @(145, 146),
    @(155, 156),
    @(159, 160),
    @(160, 322)
)

$r = 2
foreach ($row in $sheet1Data) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# ----- Sheet2 -----
$ws2 = $wb.Worksheets.Item("Sheet2")

# Clear old data range (A2:B9) before writing the new, longer table
$ws2.Range("A2:B9").ClearContents()

$sheet2Data = @(
    @(1, 8),
    @(2, 3),
    @(4, 10),
    @(5, 6),
    @(90, 91),
    @(140, 141),
    @(152, 153),
    @(154, 155),
    @(159, 160),
    @(163, 164),
    @(181, 182)
)

$r = 2
foreach ($row in $sheet2Data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r++
}
